$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Range("H39").Value = 404.94116
$ws.Range("I39").Value = 125.77778
$ws.Range("J39").Value = 719
$ws.Range("K39").Value = 377.33334
$ws.Range("L39").Value = 2157
$ws.Range("M39").Value = -81.33334000000002
$ws.Range("N39").Value = -2749
# Row 75
$ws.Range("H75").Value = 35000
$ws.Range("I75").Value = 10000
$ws.Range("K75").Value = 10000
$ws.Range("M75").Value = -9064
# Row 78
$ws.Range("H78").Value = 35000
$ws.Range("I78").Value = 10000
$ws.Range("K78").Value = 30000
$ws.Range("M78").Value = -25320
# Row 130
$ws.Range("H130").Value = 49000
$ws.Range("J130").Value = 49000
$ws.Range("L130").Value = 49000
$ws.Range("N130").Value = -59040
# Row 132
$ws.Range("H132").Value = 33339322
$ws.Range("I132").Value = 39479860
$ws.Range("K132").Value = 118439580
$ws.Range("M132").Value = -118437050

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 38
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
# Row 42
$ws.Range("H42").Value = 6000
$ws.Range("J42").Value = 6000
$ws.Range("L42").Value = 6000
$ws.Range("N42").Value = -6972
# Row 74
$ws.Range("H74").Value = 1590.5186
$ws.Range("I74").Value = 1750.2142
$ws.Range("K74").Value = 1750.2142
$ws.Range("M74").Value = -876.2141999999999
# Row 77
$ws.Range("H77").Value = 1590.5186
$ws.Range("I77").Value = 1750.2142
$ws.Range("K77").Value = 8751.071
$ws.Range("M77").Value = -4383.071
# Row 132
$ws.Range("H132").Value = 1681582.4
$ws.Range("I132").Value = 797.5
$ws.Range("J132").Value = 8404722
$ws.Range("K132").Value = 2392.5
$ws.Range("L132").Value = 25214166
$ws.Range("M132").Value = 137.5
$ws.Range("N132").Value = -25219226
# Row 133
$ws.Range("H133").Value = 39500
$ws.Range("J133").Value = 39500
$ws.Range("L133").Value = 39500
$ws.Range("N133").Value = -44560

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 38
$ws.Range("H38").Value = 30000
$ws.Range("J38").Value = 30000
$ws.Range("L38").Value = 30000
$ws.Range("N38").Value = -30832
# Row 86
$ws.Range("H86").Value = 2718.182
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 3580
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 3580
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -5826
# Row 89
$ws.Range("H89").Value = 2718.182
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 3580
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 17900
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -29132
# Row 134
$ws.Range("H134").Value = 2367371.8
$ws.Range("I134").Value = 1121.0731
$ws.Range("K134").Value = 3363.2193
$ws.Range("M134").Value = -828.2193000000002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1389.0476
$ws.Range("I16").Value = 1075.3846
$ws.Range("J16").Value = 1898.75
$ws.Range("K16").Value = 1075.3846
$ws.Range("L16").Value = 1898.75
$ws.Range("M16").Value = -788.3846000000001
$ws.Range("N16").Value = -2472.75
# Row 35
$ws.Range("H35").Value = 1667399.6
$ws.Range("I35").Value = 1667399.6
$ws.Range("K35").Value = 1667399.6
$ws.Range("M35").Value = -1667105.6
# Row 38
$ws.Range("H38").Value = 11805.5
$ws.Range("I38").Value = 3799
$ws.Range("J38").Value = 12276.471
$ws.Range("K38").Value = 3799
$ws.Range("L38").Value = 12276.471
$ws.Range("M38").Value = -3422
$ws.Range("N38").Value = -13030.471
# Row 46
$ws.Range("H46").Value = 11805.5
$ws.Range("I46").Value = 3799
$ws.Range("J46").Value = 12276.471
$ws.Range("K46").Value = 3799
$ws.Range("L46").Value = 12276.471
$ws.Range("M46").Value = -3588
$ws.Range("N46").Value = -12698.471
# Row 99
$ws.Range("H99").Value = 50002320
$ws.Range("I99").Value = 71430660
$ws.Range("J99").Value = 2866.6667
$ws.Range("K99").Value = 71430660
$ws.Range("L99").Value = 2866.6667
$ws.Range("M99").Value = -71429162
$ws.Range("N99").Value = -5862.6667
# Row 105
$ws.Range("H105").Value = 2604.403
$ws.Range("I105").Value = 2790.0908
$ws.Range("J105").Value = 1753.3334
$ws.Range("K105").Value = 2790.0908
$ws.Range("L105").Value = 1753.3334
$ws.Range("M105").Value = -1043.0908
$ws.Range("N105").Value = -5247.3334
# Row 107
$ws.Range("H107").Value = 433.80646
$ws.Range("I107").Value = 326.72726
$ws.Range("J107").Value = 695.55554
$ws.Range("K107").Value = 326.72726
$ws.Range("L107").Value = 695.55554
$ws.Range("M107").Value = 1593.27274
$ws.Range("N107").Value = -4535.55554
# Row 113
$ws.Range("H113").Value = 1389.0476
$ws.Range("I113").Value = 1075.3846
$ws.Range("J113").Value = 1898.75
$ws.Range("K113").Value = 1075.3846
$ws.Range("L113").Value = 1898.75
$ws.Range("M113").Value = 1094.6154
$ws.Range("N113").Value = -6238.75
# Row 126
$ws.Range("H126").Value = 50002320
$ws.Range("I126").Value = 71430660
$ws.Range("J126").Value = 2866.6667
$ws.Range("K126").Value = 214291980
$ws.Range("L126").Value = 8600.000100000001
$ws.Range("M126").Value = -214289510
$ws.Range("N126").Value = -13540.0001
# Row 132
$ws.Range("H132").Value = 27779980
$ws.Range("I132").Value = 1813.875
$ws.Range("J132").Value = 83336310
$ws.Range("K132").Value = 5441.625
$ws.Range("L132").Value = 250008930
$ws.Range("M132").Value = -2911.625
$ws.Range("N132").Value = -250013990
# Row 134
$ws.Range("H134").Value = 899.4186
$ws.Range("I134").Value = 808.34283
$ws.Range("J134").Value = 1297.875
$ws.Range("K134").Value = 2425.02849
$ws.Range("L134").Value = 3893.625
$ws.Range("M134").Value = 109.9715099999999
$ws.Range("N134").Value = -8963.625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 129
$ws.Range("H129").Value = 1521.1875
$ws.Range("I129").Value = 1008
$ws.Range("J129").Value = 1754.4546
$ws.Range("K129").Value = 3024
$ws.Range("L129").Value = 5263.3638
$ws.Range("M129").Value = 1976
$ws.Range("N129").Value = -15263.3638
# Row 131
$ws.Range("H131").Value = 777.6767599999999
$ws.Range("J131").Value = 780.5102000000001
$ws.Range("L131").Value = 2341.5306
$ws.Range("N131").Value = -12421.5306
# Row 133
$ws.Range("H133").Value = 17545246
$ws.Range("I133").Value = 19608846
$ws.Range("J133").Value = 4650
$ws.Range("K133").Value = 58826538
$ws.Range("L133").Value = 13950
$ws.Range("M133").Value = -58821478
$ws.Range("N133").Value = -24070

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 20838860
$ws.Range("I122").Value = 26322254
$ws.Range("J122").Value = 1960
$ws.Range("K122").Value = 78966762
$ws.Range("L122").Value = 5880
$ws.Range("M122").Value = -78964312
$ws.Range("N122").Value = -10780
# Row 126
$ws.Range("H126").Value = 2366.6667
$ws.Range("I126").Value = 1700
$ws.Range("J126").Value = 2588.889
$ws.Range("K126").Value = 5100
$ws.Range("L126").Value = 7766.667
$ws.Range("M126").Value = -2630
$ws.Range("N126").Value = -12706.667
# Row 132
$ws.Range("H132").Value = 9147.5
$ws.Range("I132").Value = 1306.7
$ws.Range("K132").Value = 3920.1
$ws.Range("M132").Value = -1390.1

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 133
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 14851.667
$ws.Range("J75").Value = 18082.5
$ws.Range("L75").Value = 18082.5
$ws.Range("N75").Value = -19954.5
# Row 78
$ws.Range("H78").Value = 14851.667
$ws.Range("J78").Value = 18082.5
$ws.Range("L78").Value = 54247.5
$ws.Range("N78").Value = -63607.5
# Row 132
$ws.Range("H132").Value = 29900.334
$ws.Range("I132").Value = 73438.07000000001
$ws.Range("J132").Value = 5519.2
$ws.Range("K132").Value = 220314.21
$ws.Range("L132").Value = 16557.6
$ws.Range("M132").Value = -217784.21
$ws.Range("N132").Value = -21617.6
# Row 136
$ws.Range("H136").Value = 1725.8959
$ws.Range("I136").Value = 1307.3889
$ws.Range("J136").Value = 1977
$ws.Range("K136").Value = 3922.1667
$ws.Range("L136").Value = 5931
$ws.Range("M136").Value = -1372.1667
$ws.Range("N136").Value = -11031

